$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.039.98'
$ws.Range("E2").Value = '  +3.49%  '

$ws.Range("D3").Value = '2.499.86'
$ws.Range("E3").Value = '  +1.75%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '494.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.21%  '

$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.515'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.86%  '

$ws.Range("D9").Value = '2.514.15'
$ws.Range("E9").Value = '  +2.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.100'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.75'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.337'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.96%  '

$ws.Range("E13").Value = '  +1.33%  '

$ws.Range("D14").Value = '2.936.96'
$ws.Range("E14").Value = '  +2.63%  '

$ws.Range("D15").Value = '57.239.63'
$ws.Range("E15").Value = '  +3.70%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.94%  '

$ws.Range("E17").Value = '  +2.94%  '

$ws.Range("D18").Value = '2.519.58'
$ws.Range("E18").Value = '  +2.94%  '

$ws.Range("E19").Value = '  +5.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.71%  '

$ws.Range("E22").Value = '  +0.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.83%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.410'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.72%  '

$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.22%  '

$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.163'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.04%  '

$ws.Range("D28").Value = '2.619.31'
$ws.Range("E28").Value = '  +2.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.15%  '

$ws.Range("D30").Value = '0.0₃0824'
$ws.Range("E30").Value = '  +6.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '151.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.54%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.53'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.32'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.22%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.27'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.47%  '

$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.81'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.04%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.15'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.49%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.885'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.93%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '34.41'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.53'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.618'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.99%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0560'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.53%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.995'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.50%  '

$ws.Range("E45").Value = '  +6.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '268.14'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.60%  '

$ws.Range("E47").Value = '  +5.36%  '

$ws.Range("E48").Value = '  +3.62%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.21'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.90%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.47%  '

$ws.Range("D51").Value = '1.892.31'
$ws.Range("E51").Value = '  -2.04%  '
